$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Ligand/Receptor/Edge expression statistics recomputed with new TPM values
$ws.Range("G2").Value = 13.16594766666667
$ws.Range("H2").Value = 39.497843
$ws.Range("I2").Value = 0.6940777873489595
$ws.Range("J2").Value = 0.6940777873489595
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.081040666666667
$ws.Range("N2").Value = 24.243122
$ws.Range("O2").Value = 0.4661250698616886
$ws.Range("P2").Value = 0.4661250698616886
$ws.Range("Q2").Value = 106.3945585095385
$ws.Range("R2").Value = 957.5510265858461
$ws.Range("S2").Value = 0.32352705711748
$ws.Range("T2").Value = 0.32352705711748
$ws.Range("G3").Value = 13.16594766666667
$ws.Range("H3").Value = 39.497843
$ws.Range("I3").Value = 0.6940777873489595
$ws.Range("J3").Value = 0.6940777873489595
$ws.Range("O3").Value = 0.4037865631294714
$ws.Range("P3").Value = 0.4037865631294715
$ws.Range("Q3").Value = 92.16559222826514
$ws.Range("R3").Value = 829.4903300543862
$ws.Range("S3").Value = 0.2802592842981445
$ws.Range("T3").Value = 0.2802592842981445
$ws.Range("G4").Value = 13.16594766666667
$ws.Range("H4").Value = 39.497843
$ws.Range("I4").Value = 0.6940777873489595
$ws.Range("J4").Value = 0.6940777873489595
$ws.Range("O4").Value = 0.1300883670088399
$ws.Range("P4").Value = 0.1300883670088399
$ws.Range("Q4").Value = 29.69309155424578
$ws.Range("R4").Value = 267.237823988212
$ws.Range("S4").Value = 0.09029144593333499
$ws.Range("T4").Value = 0.09029144593333499
$ws.Range("I5").Value = 0.1706596770095176
$ws.Range("J5").Value = 0.1706596770095176
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.081040666666667
$ws.Range("N5").Value = 24.243122
$ws.Range("O5").Value = 0.4661250698616886
$ws.Range("P5").Value = 0.4661250698616886
$ws.Range("Q5").Value = 26.16026808776
$ws.Range("R5").Value = 235.44241278984
$ws.Range("S5").Value = 0.07954875386863461
$ws.Range("T5").Value = 0.07954875386863461
$ws.Range("I6").Value = 0.1706596770095176
$ws.Range("J6").Value = 0.1706596770095176
$ws.Range("O6").Value = 0.4037865631294714
$ws.Range("P6").Value = 0.4037865631294715
$ws.Range("S6").Value = 0.06891008444445879
$ws.Range("T6").Value = 0.0689100844444588
$ws.Range("I7").Value = 0.1706596770095176
$ws.Range("J7").Value = 0.1706596770095176
$ws.Range("O7").Value = 0.1300883670088399
$ws.Range("P7").Value = 0.1300883670088399
$ws.Range("S7").Value = 0.02220083869642421
$ws.Range("T7").Value = 0.02220083869642421
$ws.Range("H8").Value = 7.697376999999999
$ws.Range("I8").Value = 0.1352625356415228
$ws.Range("J8").Value = 0.1352625356415228
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 8.081040666666667
$ws.Range("N8").Value = 24.243122
$ws.Range("O8").Value = 0.4661250698616886
$ws.Range("P8").Value = 0.4661250698616886
$ws.Range("Q8").Value = 20.73427218788822
$ws.Range("R8").Value = 186.608449690994
$ws.Range("S8").Value = 0.06304925887557394
$ws.Range("T8").Value = 0.06304925887557396
$ws.Range("H9").Value = 7.697376999999999
$ws.Range("I9").Value = 0.1352625356415228
$ws.Range("J9").Value = 0.1352625356415228
$ws.Range("O9").Value = 0.4037865631294714
$ws.Range("P9").Value = 0.4037865631294715
$ws.Range("S9").Value = 0.05461719438686811
$ws.Range("T9").Value = 0.05461719438686813
$ws.Range("H10").Value = 7.697376999999999
$ws.Range("I10").Value = 0.1352625356415228
$ws.Range("J10").Value = 0.1352625356415228
$ws.Range("O10").Value = 0.1300883670088399
$ws.Range("P10").Value = 0.1300883670088399
$ws.Range("Q10").Value = 5.786617765140887
$ws.Range("R10").Value = 52.07955988626799
$ws.Range("S10").Value = 0.0175960823790807
$ws.Range("T10").Value = 0.01759608237908071

Write-Output "Updated Efnb1-Epha4 LR-pair TPM statistics"
